$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "30.435.79"
$ws.Range("E2").Value = "  +0.40%  "

$ws.Range("D3").Value = "1.864.97"
$ws.Range("E3").Value = "  -0.16%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.001"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.08%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "235.53"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.77%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.001"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.05%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4810"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.12%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2790"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.27%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06552"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.91%  "

$ws.Range("D10").Value = "1.913.31"
$ws.Range("E10").Value = "  +2.51%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07446"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.00%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "16.16"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -1.66%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.076"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +0.29%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "86.89"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -1.01%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.6380"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -2.21%  "

$ws.Range("D16").Value = "30.416.06"
$ws.Range("E16").Value = "  +0.48%  "

$ws.Range("E17").Value = "  -0.08%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "12.96"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -1.99%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "230.19"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +5.19%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.000007460"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -1.54%  "

$ws.Range("D21").Value = "2.111.89"
$ws.Range("E21").Value = "  +0.19%  "

$ws.Range("E22").Value = "  -0.13%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.135"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -2.61%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "6.090"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.83%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "168.76"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.29%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "9.319"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.34%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "18.14"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -1.35%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.895"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -3.56%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.1055"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +12.76%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.385"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -4.63%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.243"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -1.14%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.968"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.98%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.04965"
$ws.Range("D33").Style = "Normal"

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.169"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -2.51%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.7391"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -0.71%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.9995"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.19%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.711"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.01%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.01932"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +6.31%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.636"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.86%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.9150"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.90%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.033"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -1.88%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "105.83"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.55%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.9959"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.70%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.4169"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -1.90%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "5.549"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -6.14%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "7.158"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -2.79%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "61.85"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -2.81%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.1222"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -4.76%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "8.891"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.13%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.417"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -3.81%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "33.42"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.33%  "
